# Age Photo Exp Finished
# Applies the "Exp 2 Alignment Comparison" rework on Sheet3, adding the
# l0/l1/l2 alignment rows and the new 2.1/2.2 alignment columns (I:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Clear out the old "Exp 2 Alignment Comparison" block (rows 15-27); it will
# be fully rebuilt below with the new layout (rows 15-24 and 28-31).
$ws.Range("A15:L27").ClearContents()

# --- Row 15 / 16: create the new shared strings "Align 2.1" / "Align 2.2"
# first (columns I), so they land at the expected shared-string indices,
# then fill in the rest of row 15/16's data.
$ws.Cells.Item(15, 9).Value = "Align 2.1"
$ws.Cells.Item(15, 10).Value = 5.138
$ws.Cells.Item(15, 11).Value = 9.5967
$ws.Cells.Item(15, 12).Value = 12.1284

$ws.Cells.Item(16, 9).Value = "Align 2.2"
$ws.Cells.Item(16, 10).Value = 4.1008
$ws.Cells.Item(16, 11).Value = 7.0898
$ws.Cells.Item(16, 12).Value = 10.3975

# --- Now create "Align l2" / "Align l1" / "Align l0" in that order (A17,
# A16, A15) so the new shared strings come out in the same order as the
# target workbook.
$ws.Cells.Item(17, 1).Value = "Align l2"
$ws.Cells.Item(17, 2).Value = 3.864211
$ws.Cells.Item(17, 3).Value = 6.756754
$ws.Cells.Item(17, 4).Value = 9.926225

$ws.Cells.Item(16, 1).Value = "Align l1"
$ws.Cells.Item(16, 2).Value = 4.096274
$ws.Cells.Item(16, 3).Value = 7.060792
$ws.Cells.Item(16, 4).Value = 10.179871

$ws.Cells.Item(15, 1).Value = "Align l0"
$ws.Cells.Item(15, 2).Value = 5.004651
$ws.Cells.Item(15, 3).Value = 11.678655
$ws.Cells.Item(15, 4).Value = 11.845544

# --- Row 19 / 20: second mini-table for the 2.1 / 2.2 columns.
$ws.Cells.Item(19, 9).Value = "Align 2.1"
$ws.Cells.Item(19, 10).Value = 0.5782
$ws.Cells.Item(19, 11).Value = 0.3436
$ws.Cells.Item(19, 12).Value = 0.3106

$ws.Cells.Item(20, 9).Value = "Align 2.2"

# --- Row 21: "Accuracy within 5 year error" header (moved down from 19).
$ws.Cells.Item(21, 1).Value = "Accuracy within 5 year error"

# --- Rows 22-24: l0/l1/l2 accuracy rows.
$ws.Cells.Item(22, 1).Value = "Align l0"
$ws.Cells.Item(22, 2).Value = 0.584416
$ws.Cells.Item(22, 3).Value = 0.289195
$ws.Cells.Item(22, 4).Value = 0.31924

$ws.Cells.Item(23, 1).Value = "Align l1"
$ws.Cells.Item(23, 2).Value = 0.680064
$ws.Cells.Item(23, 3).Value = 0.447004
$ws.Cells.Item(23, 4).Value = 0.366426

$ws.Cells.Item(24, 1).Value = "Align l2"
$ws.Cells.Item(24, 2).Value = 0.710964
$ws.Cells.Item(24, 3).Value = 0.473805
$ws.Cells.Item(24, 4).Value = 0.374806

# --- Rows 28-31: "Accuracy of Gender" block (moved down from rows 24-27).
$ws.Cells.Item(28, 1).Value = "Accuracy of Gender"
$ws.Cells.Item(29, 1).Value = "Align 0"
$ws.Cells.Item(30, 1).Value = "Align 1"
$ws.Cells.Item(31, 1).Value = "Align 2"

# --- Update the view: scroll so row 10 is at the top and select A26, to
# match the author's final cursor position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A26").Select()
